$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.178.19'
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("D3").Value = '3.403.75'
$ws.Range("E3").Value = '  -4.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.91'
$ws.Range("E5").Value = '  -4.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.43'
$ws.Range("E6").Value = '  -9.07%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.400.98'
$ws.Range("E8").Value = '  -4.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.480'
$ws.Range("E9").Value = '  -7.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.120'
$ws.Range("E10").Value = '  -9.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.95'
$ws.Range("E11").Value = '  -10.71%  '
$ws.Range("E12").Value = '  -9.99%  '
$ws.Range("D13").Value = '3.977.38'
$ws.Range("E13").Value = '  -4.27%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '3.445.87'
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000176'
$ws.Range("E15").Value = '  -10.47%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.115'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.95'
$ws.Range("E17").Value = '  -10.53%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '64.225.48'
$ws.Range("E18").Value = '  -3.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.50'
$ws.Range("E19").Value = '  -13.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.64'
$ws.Range("E20").Value = '  -9.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.53'
$ws.Range("E21").Value = '  -7.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '377.82'
$ws.Range("E22").Value = '  -11.48%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.540'
$ws.Range("E24").Value = '  -9.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.71'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.27'
$ws.Range("E26").Value = '  -8.07%  '
$ws.Range("D27").Value = '3.537.46'
$ws.Range("E27").Value = '  -4.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000103'
$ws.Range("E28").Value = '  -10.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.09'
$ws.Range("E30").Value = '  -11.66%  '
$ws.Range("E31").Value = '  -12.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.97'
$ws.Range("E32").Value = '  -11.86%  '
$ws.Range("D33").Value = '3.414.37'
$ws.Range("E33").Value = '  -4.05%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.86'
$ws.Range("E35").Value = '  -6.63%  '
$ws.Range("E36").Value = '  -11.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '170.81'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.17'
$ws.Range("E38").Value = '  -13.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.63'
$ws.Range("E39").Value = '  -13.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.44'
$ws.Range("E40").Value = '  -12.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.59'
$ws.Range("E41").Value = '  -13.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0754'
$ws.Range("E42").Value = '  -8.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.794'
$ws.Range("E43").Value = '  -8.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.57'
$ws.Range("E45").Value = '  -8.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.25'
$ws.Range("E46").Value = '  -15.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.59'
$ws.Range("E47").Value = '  -11.08%  '
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.07'
$ws.Range("E49").Value = '  -5.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.48'
$ws.Range("E50").Value = '  -9.21%  '
$ws.Range("D51").Value = '2.181.73'
$ws.Range("E51").Value = '  -6.74%  '
